# Auto-generated Excel COM-interop script to apply the market-data update
# described by the commit "chore: update Sheets via scheduled runner".
# For each affected leve row, columns H:N (currentAveragePrice, currentAveragePriceNQ,
# currentAveragePriceHQ, LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ) are refreshed
# with newly pulled market values. Where a column has no applicable value it is cleared
# (cell removed) and where a new value becomes applicable a cell is written.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 76 (Leve Item ID 12602)
$ws.Range("H76").Value = 52990.15
$ws.Range("I76").Value = 52990.15
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 52990.15
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -52675.15
$ws.Range("N76").ClearContents()

# ALC row 79 (Leve Item ID 12602)
$ws.Range("H79").Value = 52990.15
$ws.Range("I79").Value = 52990.15
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 52990.15
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -51898.15
$ws.Range("N79").ClearContents()

# ALC row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 7484.5
$ws.Range("I98").Value = 7689.4707
$ws.Range("K98").Value = 7689.4707
$ws.Range("M98").Value = -6191.4707

# ALC row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 7484.5
$ws.Range("I122").Value = 7689.4707
$ws.Range("K122").Value = 23068.4121
$ws.Range("M122").Value = -20618.4121

# ALC row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 4205219
$ws.Range("I132").Value = 5294594
$ws.Range("K132").Value = 15883782
$ws.Range("M132").Value = -15881252

# ALC row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1132.2778
$ws.Range("I137").Value = 1059.174
$ws.Range("K137").Value = 3177.522
$ws.Range("M137").Value = -627.5219999999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 14389.988
$ws.Range("I32").Value = 16249.945
$ws.Range("J32").Value = 3075.25
$ws.Range("K32").Value = 16249.945
$ws.Range("L32").Value = 3075.25
$ws.Range("M32").Value = -15962.945
$ws.Range("N32").Value = -3649.25

# ARM row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1546.55
$ws.Range("I61").Value = 1239.0667
$ws.Range("J61").Value = 2469
$ws.Range("K61").Value = 1239.0667
$ws.Range("L61").Value = 2469
$ws.Range("M61").Value = -1027.0667
$ws.Range("N61").Value = -2893

# ARM row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1089.3334
$ws.Range("I74").Value = 1104.4
$ws.Range("J74").Value = 1014
$ws.Range("K74").Value = 1104.4
$ws.Range("L74").Value = 1014
$ws.Range("M74").Value = -230.4000000000001
$ws.Range("N74").Value = -2762

# ARM row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1089.3334
$ws.Range("I77").Value = 1104.4
$ws.Range("J77").Value = 1014
$ws.Range("K77").Value = 5522
$ws.Range("L77").Value = 5070
$ws.Range("M77").Value = -1154
$ws.Range("N77").Value = -13806

# ARM row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 6266.276
$ws.Range("I132").Value = 7674.6875
$ws.Range("J132").Value = 4532.846
$ws.Range("K132").Value = 23024.0625
$ws.Range("L132").Value = 13598.538
$ws.Range("M132").Value = -20494.0625
$ws.Range("N132").Value = -18658.538

# ARM row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1546.55
$ws.Range("I136").Value = 1239.0667
$ws.Range("J136").Value = 2469
$ws.Range("K136").Value = 3717.2001
$ws.Range("L136").Value = 7407
$ws.Range("M136").Value = -1167.2001
$ws.Range("N136").Value = -12507

$ws = $wb.Worksheets.Item("BSM")
# BSM row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 46150.043
$ws.Range("I134").Value = 55234.26
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 165702.78
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -163167.78
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3877395
$ws.Range("I31").Value = 1844.0333
$ws.Range("J31").Value = 5953583
$ws.Range("K31").Value = 1844.0333
$ws.Range("L31").Value = 5953583
$ws.Range("M31").Value = -1549.0333
$ws.Range("N31").Value = -5954173

# CRP row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3877395
$ws.Range("I34").Value = 1844.0333
$ws.Range("J34").Value = 5953583
$ws.Range("K34").Value = 1844.0333
$ws.Range("L34").Value = 5953583
$ws.Range("M34").Value = -1642.0333
$ws.Range("N34").Value = -5953987

# CRP row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 2666.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12900.0001

# CRP row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2249.0732
$ws.Range("I132").Value = 1857.6428
$ws.Range("J132").Value = 3092.1538
$ws.Range("K132").Value = 5572.928400000001
$ws.Range("L132").Value = 9276.4614
$ws.Range("M132").Value = -3042.928400000001
$ws.Range("N132").Value = -14336.4614

# CRP row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1000.86664
$ws.Range("I134").Value = 1011.6429
$ws.Range("K134").Value = 3034.9287
$ws.Range("M134").Value = -499.9287000000004

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1390.1957
$ws.Range("I68").Value = 1203.7059
$ws.Range("J68").Value = 1918.5834
$ws.Range("K68").Value = 3611.1177
$ws.Range("L68").Value = 5755.7502
$ws.Range("M68").Value = -2800.1177
$ws.Range("N68").Value = -7377.7502

# CUL row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1390.1957
$ws.Range("I71").Value = 1203.7059
$ws.Range("J71").Value = 1918.5834
$ws.Range("K71").Value = 10833.3531
$ws.Range("L71").Value = 17267.2506
$ws.Range("M71").Value = -6777.3531
$ws.Range("N71").Value = -25379.2506

# CUL row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 662.9530999999999
$ws.Range("I113").Value = 1090.5834
$ws.Range("J113").Value = 564.2692
$ws.Range("K113").Value = 3271.7502
$ws.Range("L113").Value = 1692.8076
$ws.Range("M113").Value = -1101.7502
$ws.Range("N113").Value = -6032.8076

# CUL row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 882
$ws.Range("J122").Value = 1696.9166
$ws.Range("L122").Value = 15272.2494
$ws.Range("N122").Value = -20172.2494

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 4820990.5
$ws.Range("I131").Value = 8680
$ws.Range("J131").Value = 8217915.5
$ws.Range("K131").Value = 26040
$ws.Range("L131").Value = 24653746.5
$ws.Range("M131").Value = -21000
$ws.Range("N131").Value = -24663826.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 7850694
$ws.Range("I70").Value = 10412370
$ws.Range("J70").Value = 5562.4375
$ws.Range("K70").Value = 10412370
$ws.Range("L70").Value = 5562.4375
$ws.Range("M70").Value = -10412100
$ws.Range("N70").Value = -6102.4375

# GSM row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 7850694
$ws.Range("I73").Value = 10412370
$ws.Range("J73").Value = 5562.4375
$ws.Range("K73").Value = 10412370
$ws.Range("L73").Value = 5562.4375
$ws.Range("M73").Value = -10411434
$ws.Range("N73").Value = -7434.4375

# GSM row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 6026.6113
$ws.Range("I80").Value = 4675.385
$ws.Range("J80").Value = 9539.799999999999
$ws.Range("K80").Value = 4675.385
$ws.Range("L80").Value = 9539.799999999999
$ws.Range("M80").Value = -3677.385
$ws.Range("N80").Value = -11535.8

# GSM row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 6026.6113
$ws.Range("I83").Value = 4675.385
$ws.Range("J83").Value = 9539.799999999999
$ws.Range("K83").Value = 23376.925
$ws.Range("L83").Value = 47699
$ws.Range("M83").Value = -18384.925
$ws.Range("N83").Value = -57683

# GSM row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 50004900
$ws.Range("I122").Value = 100006800
$ws.Range("K122").Value = 300020400
$ws.Range("M122").Value = -300017950

# GSM row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 13200
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 25000
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 75000
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -79940

# GSM row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 82548.39999999999
$ws.Range("I132").Value = 107553.266
$ws.Range("J132").Value = 3366.3333
$ws.Range("K132").Value = 322659.798
$ws.Range("L132").Value = 10098.9999
$ws.Range("M132").Value = -320129.798
$ws.Range("N132").Value = -15158.9999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 2002.7646
$ws.Range("I40").Value = 2048.5
$ws.Range("J40").Value = 1962.1111
$ws.Range("K40").Value = 2048.5
$ws.Range("L40").Value = 1962.1111
$ws.Range("M40").Value = -1912.5
$ws.Range("N40").Value = -2234.1111

# LTW row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 2490.3333
$ws.Range("I122").Value = 2150.1667
$ws.Range("J122").Value = 2943.889
$ws.Range("K122").Value = 6450.500100000001
$ws.Range("L122").Value = 8831.667000000001
$ws.Range("M122").Value = -4000.500100000001
$ws.Range("N122").Value = -13731.667

# LTW row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 9530.214
$ws.Range("I132").Value = 14884.6
$ws.Range("J132").Value = 3352.077
$ws.Range("K132").Value = 44653.8
$ws.Range("L132").Value = 10056.231
$ws.Range("M132").Value = -42123.8
$ws.Range("N132").Value = -15116.231

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 1273.341
$ws.Range("I132").Value = 1166.7646
$ws.Range("J132").Value = 1635.7
$ws.Range("K132").Value = 3500.2938
$ws.Range("L132").Value = 4907.1
$ws.Range("M132").Value = -970.2937999999999
$ws.Range("N132").Value = -9967.1

